# Applies the cryptos-list price/volume refresh described in the commit
# "Updated cryptos list on Tue Jul 30 05:59:45 UTC 2024 with GitHub Actions".
# All touched cells are plain text in the original workbook (t="inlineStr"),
# including values that *look* numeric ("1.00", "73.49", ...). Excel would
# normally auto-convert a bare numeric-looking string to a real number, so
# any such value is written with a leading apostrophe (Excel's force-text
# marker) and the resulting "quote prefix" cell style is stripped back to
# Normal afterwards so the cell keeps its original (style-less) text form.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.531.32'
$ws.Range("E2").Value = '  -4.52%  '
$ws.Range("D3").Value = '3.310.10'
$ws.Range("E3").Value = '  -1.39%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = "'571.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.51%  '
$ws.Range("D6").Value = "'182.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.88%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  -1.64%  '
$ws.Range("E9").Value = '  -3.87%  '
$ws.Range("E10").Value = '  -2.30%  '
$ws.Range("E11").Value = '  -4.99%  '
$ws.Range("D12").Value = '3.885.78'
$ws.Range("E12").Value = '  -1.35%  '
$ws.Range("E13").Value = '  -0.83%  '
$ws.Range("E14").Value = '  -4.75%  '
$ws.Range("D15").Value = '66.616.54'
$ws.Range("D17").Value = '3.309.46'
$ws.Range("E17").Value = '  -1.00%  '
$ws.Range("E18").Value = '  -0.60%  '
$ws.Range("E19").Value = '  -2.66%  '
$ws.Range("D20").Value = "'431.69"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.92%  '
$ws.Range("E21").Value = '  -2.28%  '
$ws.Range("D22").Value = "'73.49"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.10%  '
$ws.Range("D23").Value = "'1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.02%  '
$ws.Range("D24").Value = "'0.519"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.54%  '
$ws.Range("E25").Value = '  -3.24%  '
$ws.Range("E26").Value = '  +0.10%  '
$ws.Range("E27").Value = '  -5.53%  '
$ws.Range("D28").Value = "'1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.36%  '
$ws.Range("E29").Value = '  -2.37%  '
$ws.Range("E30").Value = '  -1.75%  '
$ws.Range("E31").Value = '  -5.43%  '
$ws.Range("E33").Value = '  -3.90%  '
$ws.Range("E34").Value = '  -3.85%  '
$ws.Range("E35").Value = '  -1.38%  '
$ws.Range("D36").Value = "'159.58"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.07%  '
$ws.Range("E37").Value = '  -4.77%  '
$ws.Range("D38").Value = "'27.13"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.44%  '
$ws.Range("D39").Value = '2.810.29'
$ws.Range("E39").Value = '  +1.99%  '
$ws.Range("E40").Value = '  -3.24%  '
$ws.Range("E41").Value = '  -3.72%  '
$ws.Range("D42").Value = "'6.18"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.79%  '
$ws.Range("E43").Value = '  -2.06%  '
$ws.Range("D44").Value = "'40.12"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.57%  '
$ws.Range("D45").Value = "'24.38"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.37%  '
$ws.Range("E46").Value = '  -6.89%  '
$ws.Range("D47").Value = "'320.11"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -7.51%  '
$ws.Range("E48").Value = '  -4.33%  '
$ws.Range("E49").Value = '  -2.42%  '
$ws.Range("D50").Value = "'6.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.10%  '
$ws.Range("E51").Value = '  -1.26%  '
